# "quick data cleaning; added missing negatives"
#
# The Sept 2020 finance ledger had an erroneous entry on row 18 (an "Active
# Living" / COVID-refund line that was the only *positive* amount in the
# whole COST column) and a stray, empty, red-highlighted placeholder cell
# in F7 left over from an earlier review pass. Both get cleaned up:
#
#   1. Delete row 18 entirely (shifts the real last entry - "Little
#      Caesars" / Restaurant / -15, which used to be row 19 - up into
#      row 18).
#   2. Clear the leftover highlighted placeholder cell F7 (no value, just
#      a red fill) so it no longer appears in the sheet.
#   3. Leave the selection on the new last row (18) to match where the
#      editor's cursor ended up after the cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the bad "Active Living" / +382 COVID-refund row; everything
#    below shifts up so the true final row (Little Caesars, -15) becomes
#    row 18.
$ws.Rows(18).Delete()

# 2. Remove the stray empty red-filled cell.
$ws.Range("F7").Clear()

# 3. Match the post-edit selection (whole of row 18).
$ws.Rows(18).Select()
